$d = $word.ActiveDocument
# insert and delete a char right where the bookmark sits (between run2 end=49 and run3 start=49)
$r = $d.Range(49, 49)
$r.InsertBefore("X")
$r2 = $d.Range(49, 50)
Write-Output $r2.Text
$r2.Delete()
Write-Output "done"
Write-Output $d.Content.Text
